$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.3275576666666667
$ws.Range("H2").Value = 0.982673
$ws.Range("I2").Value = 0.05486041027915935
$ws.Range("J2").Value = 0.05486041027915935
$ws.Range("M2").Value = 3.390429
$ws.Range("N2").Value = 10.171287
$ws.Range("O2").Value = 0.173121426386348
$ws.Range("P2").Value = 0.173121426386348
$ws.Range("Q2").Value = 1.110561012239
$ws.Range("R2").Value = 9.995049110150999
$ws.Range("S2").Value = 0.009497512479668334
$ws.Range("T2").Value = 0.009497512479668336
$ws.Range("G3").Value = 0.3275576666666667
$ws.Range("H3").Value = 0.982673
$ws.Range("I3").Value = 0.05486041027915935
$ws.Range("J3").Value = 0.05486041027915935
$ws.Range("O3").Value = 0.5936336753560868
$ws.Range("P3").Value = 0.5936336753560868
$ws.Range("Q3").Value = 3.808115662883667
$ws.Range("R3").Value = 34.273040965953
$ws.Range("S3").Value = 0.03256698698556021
$ws.Range("T3").Value = 0.03256698698556021
$ws.Range("G4").Value = 0.3275576666666667
$ws.Range("H4").Value = 0.982673
$ws.Range("I4").Value = 0.05486041027915935
$ws.Range("J4").Value = 0.05486041027915935
$ws.Range("M4").Value = 4.546141666666667
$ws.Range("N4").Value = 13.638425
$ws.Range("O4").Value = 0.2321342018628743
$ws.Range("P4").Value = 0.2321342018628743
$ws.Range("Q4").Value = 1.489123556669445
$ws.Range("R4").Value = 13.402112010025
$ws.Range("S4").Value = 0.01273497755402248
$ws.Range("T4").Value = 0.01273497755402248
$ws.Range("G5").Value = 0.3275576666666667
$ws.Range("H5").Value = 0.982673
$ws.Range("I5").Value = 0.05486041027915935
$ws.Range("J5").Value = 0.05486041027915935
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.021752
$ws.Range("N5").Value = 0.06525600000000001
$ws.Range("O5").Value = 0.001110696394691009
$ws.Range("P5").Value = 0.001110696394691009
$ws.Range("Q5").Value = 0.007125034365333335
$ws.Range("R5").Value = 0.064125309288
$ws.Range("S5").Value = [double]"6.093325990833186E-05"
$ws.Range("T5").Value = [double]"6.093325990833186E-05"
$ws.Range("I6").Value = 0.8684635977749966
$ws.Range("J6").Value = 0.8684635977749967
$ws.Range("M6").Value = 3.390429
$ws.Range("N6").Value = 10.171287
$ws.Range("O6").Value = 0.173121426386348
$ws.Range("P6").Value = 0.173121426386348
$ws.Range("Q6").Value = 17.580652556733
$ws.Range("R6").Value = 158.225873010597
$ws.Range("S6").Value = 0.150349656811427
$ws.Range("T6").Value = 0.1503496568114271
$ws.Range("I7").Value = 0.8684635977749966
$ws.Range("J7").Value = 0.8684635977749967
$ws.Range("O7").Value = 0.5936336753560868
$ws.Range("P7").Value = 0.5936336753560868
$ws.Range("S7").Value = 0.5155492374601415
$ws.Range("T7").Value = 0.5155492374601416
$ws.Range("I8").Value = 0.8684635977749966
$ws.Range("J8").Value = 0.8684635977749967
$ws.Range("M8").Value = 4.546141666666667
$ws.Range("N8").Value = 13.638425
$ws.Range("O8").Value = 0.2321342018628743
$ws.Range("P8").Value = 0.2321342018628743
$ws.Range("Q8").Value = 23.573458437075
$ws.Range("R8").Value = 212.161125933675
$ws.Range("S8").Value = 0.2016001041164591
$ws.Range("T8").Value = 0.2016001041164592
$ws.Range("I9").Value = 0.8684635977749966
$ws.Range("J9").Value = 0.8684635977749967
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.021752
$ws.Range("N9").Value = 0.06525600000000001
$ws.Range("O9").Value = 0.001110696394691009
$ws.Range("P9").Value = 0.001110696394691009
$ws.Range("Q9").Value = 0.112792320504
$ws.Range("R9").Value = 1.015130884536
$ws.Range("S9").Value = 0.0009645993869690714
$ws.Range("T9").Value = 0.0009645993869690715
$ws.Range("G10").Value = 0.1537743333333333
$ws.Range("H10").Value = 0.461323
$ws.Range("I10").Value = 0.02575461934052592
$ws.Range("J10").Value = 0.02575461934052592
$ws.Range("M10").Value = 3.390429
$ws.Range("N10").Value = 10.171287
$ws.Range("O10").Value = 0.173121426386348
$ws.Range("P10").Value = 0.173121426386348
$ws.Range("Q10").Value = 0.5213609591890001
$ws.Range("R10").Value = 4.692248632701
$ws.Range("S10").Value = 0.004458676436269272
$ws.Range("T10").Value = 0.004458676436269273
$ws.Range("G11").Value = 0.1537743333333333
$ws.Range("H11").Value = 0.461323
$ws.Range("I11").Value = 0.02575461934052592
$ws.Range("J11").Value = 0.02575461934052592
$ws.Range("O11").Value = 0.5936336753560868
$ws.Range("P11").Value = 0.5936336753560868
$ws.Range("Q11").Value = 1.787747645400334
$ws.Range("R11").Value = 16.089728808603
$ws.Range("S11").Value = 0.01528880933651336
$ws.Range("T11").Value = 0.01528880933651336
$ws.Range("G12").Value = 0.1537743333333333
$ws.Range("H12").Value = 0.461323
$ws.Range("I12").Value = 0.02575461934052592
$ws.Range("J12").Value = 0.02575461934052592
$ws.Range("M12").Value = 4.546141666666667
$ws.Range("N12").Value = 13.638425
$ws.Range("O12").Value = 0.2321342018628743
$ws.Range("P12").Value = 0.2321342018628743
$ws.Range("Q12").Value = 0.6990799040305556
$ws.Range("R12").Value = 6.291719136275002
$ws.Range("S12").Value = 0.005978528004895129
$ws.Range("T12").Value = 0.00597852800489513
$ws.Range("G13").Value = 0.1537743333333333
$ws.Range("H13").Value = 0.461323
$ws.Range("I13").Value = 0.02575461934052592
$ws.Range("J13").Value = 0.02575461934052592
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.021752
$ws.Range("N13").Value = 0.06525600000000001
$ws.Range("O13").Value = 0.001110696394691009
$ws.Range("P13").Value = 0.001110696394691009
$ws.Range("Q13").Value = 0.003344899298666667
$ws.Range("R13").Value = 0.03010409368800001
$ws.Range("S13").Value = [double]"2.860556284816147E-05"
$ws.Range("T13").Value = [double]"2.860556284816147E-05"
$ws.Range("G14").Value = 0.3040386666666667
$ws.Range("H14").Value = 0.9121160000000001
$ws.Range("I14").Value = 0.05092137260531806
$ws.Range("J14").Value = 0.05092137260531806
$ws.Range("M14").Value = 3.390429
$ws.Range("N14").Value = 10.171287
$ws.Range("O14").Value = 0.173121426386348
$ws.Range("P14").Value = 0.173121426386348
$ws.Range("Q14").Value = 1.030821512588
$ws.Range("R14").Value = 9.277393613292
$ws.Range("S14").Value = 0.008815580658983368
$ws.Range("T14").Value = 0.00881558065898337
$ws.Range("G15").Value = 0.3040386666666667
$ws.Range("H15").Value = 0.9121160000000001
$ws.Range("I15").Value = 0.05092137260531806
$ws.Range("J15").Value = 0.05092137260531806
$ws.Range("O15").Value = 0.5936336753560868
$ws.Range("P15").Value = 0.5936336753560868
$ws.Range("Q15").Value = 3.534688778430668
$ws.Range("R15").Value = 31.81219900587601
$ws.Range("S15").Value = 0.03022864157387171
$ws.Range("T15").Value = 0.03022864157387171
$ws.Range("G16").Value = 0.3040386666666667
$ws.Range("H16").Value = 0.9121160000000001
$ws.Range("I16").Value = 0.05092137260531806
$ws.Range("J16").Value = 0.05092137260531806
$ws.Range("M16").Value = 4.546141666666667
$ws.Range("N16").Value = 13.638425
$ws.Range("O16").Value = 0.2321342018628743
$ws.Range("P16").Value = 0.2321342018628743
$ws.Range("Q16").Value = 1.382202850811111
$ws.Range("R16").Value = 12.4398256573
$ws.Range("S16").Value = 0.01182059218749754
$ws.Range("T16").Value = 0.01182059218749754
$ws.Range("G17").Value = 0.3040386666666667
$ws.Range("H17").Value = 0.9121160000000001
$ws.Range("I17").Value = 0.05092137260531806
$ws.Range("J17").Value = 0.05092137260531806
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.021752
$ws.Range("N17").Value = 0.06525600000000001
$ws.Range("O17").Value = 0.001110696394691009
$ws.Range("P17").Value = 0.001110696394691009
$ws.Range("Q17").Value = 0.006613449077333336
$ws.Range("R17").Value = 0.05952104169600002
$ws.Range("S17").Value = [double]"5.655818496544428E-05"
$ws.Range("T17").Value = [double]"5.655818496544428E-05"
